$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.945.63"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").Value = "2.563.06"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("E13").Value = "  +7.14%  "

$ws.Range("D14").Value = "2.532.60"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "42.964.67"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "0.0₃0998"
$ws.Range("E18").Value = "  +3.40%  "

$ws.Range("E19").Value = "  +3.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.65"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "253.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "28.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.33%  "

$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("E29").Value = "  -3.64%  "

$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("E33").Value = "  -2.51%  "

$ws.Range("E34").Value = "  -5.68%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  -2.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.86%  "

$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.68%  "

$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0310"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.099.14"
$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +26.10%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.35%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.810.72"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("E51").Value = "  +3.17%  "
